# Workbook: medicamentos.xlsx
# Sheet "universal" (first sheet) holds a table of medication dosage rules.
# Column O ("formula") contains JS-style template-literal expressions used to
# compute the dose. This edit:
#   1. Fixes the buggy formula in row 2 (dipirona) from the invalid
#      "(p*0.04)+`  ml`" to a proper template literal "`${p*0.04} ml`".
#   2. Adds the missing formulas for rows 3 (pediatrica) and 4 (adulta) in
#      column O, which were previously empty.
#   3. Updates the sheet's active selection to O4 (last edited cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("O2").Value = '`${p*0.04} ml`'
$ws.Range("O3").Value = '`${p*2} a funcionar x2 `'
$ws.Range("O4").Value = '`a funcionar`'

$ws.Range("O4").Select() | Out-Null
